$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 396, pushing existing rows 396-414 down to 397-415
$ws.Rows.Item(396).Insert()

# Populate the newly inserted row 396 with the new data record
$ws.Cells.Item(396, 1).Value = 5
$ws.Cells.Item(396, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(396, 3).Value = "Maule"
$ws.Cells.Item(396, 4).Value = 45041
$ws.Cells.Item(396, 5).Value = 7
$ws.Cells.Item(396, 6).Value = 100112008
$ws.Cells.Item(396, 7).Value = "Coliflor"
$ws.Cells.Item(396, 8).Value = "Sin especificar"
$ws.Cells.Item(396, 9).Value = "Primera"
$ws.Cells.Item(396, 10).Value = 3000
$ws.Cells.Item(396, 11).Value = 800
$ws.Cells.Item(396, 12).Value = 800
$ws.Cells.Item(396, 13).Value = 800
$ws.Cells.Item(396, 14).Value = "$/unidad"
$ws.Cells.Item(396, 15).Value = "Región del Maule"
$ws.Cells.Item(396, 16).Value = 800
$ws.Cells.Item(396, 17).Value = 1
$ws.Cells.Item(396, 18).Value = "Hortaliza"
